# Updated cryptos list on Mon Sep  9 02:50:20 UTC 2024 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row,
# and re-sort two pairs of rows whose relative order changed (Aptos/PancakeSwap
# and Hedera/Bittensor). Numeric-looking price strings are written with a
# leading apostrophe so Excel keeps them as text (matching the source data,
# which stores these as plain strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.715.31'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '2.281.94'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''506.54'
$ws.Range('E5').Value = '  +2.31%  '
$ws.Range('D6').Value = '''129.13'
$ws.Range('E6').Value = '  +0.87%  '
$ws.Range('D7').Value = '''0.996'
$ws.Range('D8').Value = '''0.528'
$ws.Range('E8').Value = '  +0.33%  '
$ws.Range('D9').Value = '2.299.80'
$ws.Range('E9').Value = '  +0.78%  '
$ws.Range('D10').Value = '''0.0970'
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('D11').Value = '''0.153'
$ws.Range('E11').Value = '  +0.95%  '
$ws.Range('D12').Value = '''0.344'
$ws.Range('E12').Value = '  +4.82%  '
$ws.Range('D13').Value = '''4.91'
$ws.Range('E13').Value = '  +4.08%  '
$ws.Range('D14').Value = '''23.41'
$ws.Range('E14').Value = '  +5.80%  '
$ws.Range('D15').Value = '2.686.62'
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('D16').Value = '54.716.16'
$ws.Range('E16').Value = '  +0.97%  '
$ws.Range('E17').Value = '  +1.22%  '
$ws.Range('D18').Value = '2.307.00'
$ws.Range('E18').Value = '  +1.41%  '
$ws.Range('D19').Value = '''10.32'
$ws.Range('E19').Value = '  +2.64%  '
$ws.Range('D20').Value = '''4.14'
$ws.Range('E20').Value = '  +0.95%  '
$ws.Range('D21').Value = '''306.79'
$ws.Range('E21').Value = '  +1.34%  '
$ws.Range('D22').Value = '''6.40'
$ws.Range('E22').Value = '  -0.92%  '
$ws.Range('D23').Value = '''0.998'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '''60.19'
$ws.Range('E24').Value = '  -2.50%  '
$ws.Range('D25').Value = '''0.995'
$ws.Range('E25').Value = '  -0.37%  '
$ws.Range('D26').Value = '''0.151'
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('D27').Value = '''7.42'
$ws.Range('E27').Value = '  +3.35%  '
$ws.Range('D28').Value = '''170.49'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').Value = '0.0₃0704'
$ws.Range('E29').Value = '  +3.61%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').Value = '''6.08'
$ws.Range('E30').Value = '  +4.00%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = '''1.63'
$ws.Range('E31').Value = '  +1.38%  '
$ws.Range('E32').Value = '  +2.33%  '
$ws.Range('E33').Value = '  +0.04%  '
$ws.Range('D34').Value = '''17.97'
$ws.Range('E34').Value = '  +1.43%  '
$ws.Range('E35').Value = '  -0.18%  '
$ws.Range('D36').Value = '''0.923'
$ws.Range('E36').Value = '  +3.79%  '
$ws.Range('D37').Value = '''1.20'
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('D38').Value = '''3.79'
$ws.Range('E38').Value = '  +2.69%  '
$ws.Range('D39').Value = '''36.32'
$ws.Range('E39').Value = '  +1.56%  '
$ws.Range('D40').Value = '''0.377'
$ws.Range('E40').Value = '  +1.29%  '
$ws.Range('D41').Value = '''1.41'
$ws.Range('E41').Value = '  +0.78%  '
$ws.Range('D42').Value = '''3.41'
$ws.Range('E42').Value = '  +0.69%  '
$ws.Range('D43').Value = '''4.93'
$ws.Range('E43').Value = '  +3.56%  '
$ws.Range('D44').Value = '''125.22'
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Value = '''0.0499'
$ws.Range('E45').Value = '  +2.50%  '
$ws.Range('B46').Value = 'Bittensor'
$ws.Range('C46').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D46').Value = '''249.65'
$ws.Range('E46').Value = '  +4.80%  '
$ws.Range('D47').Value = '''0.0904'
$ws.Range('E47').Value = '  +0.78%  '
$ws.Range('D48').Value = '''0.550'
$ws.Range('E48').Value = '  +0.94%  '
$ws.Range('D49').Value = '''0.376'
$ws.Range('E49').Value = '  +1.53%  '
$ws.Range('D50').Value = '''0.0207'
$ws.Range('E50').Value = '  +0.98%  '
$ws.Range('E51').Value = '  +0.52%  '
